$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "$ 601,56"
$ws.Range("F3").Value = "$ 1.851,15"
$ws.Range("F4").Value = "$ 1.018,14"
$ws.Range("F7").Value = "$ 1.758,58"
$ws.Range("F8").Value = "$ 629,33"
$ws.Range("F10").Value = "$ 481,23"
$ws.Range("F11").Value = "$ 1.851,15"
$ws.Range("F12").Value = "$ 305,36"
$ws.Range("F13").Value = "$ 305,36"
$ws.Range("F14").Value = "$ 490,48"
$ws.Range("F15").Value = "$ 555,28"
$ws.Range("F16").Value = "$ 1.018,08"
$ws.Range("F17").Value = "$ 2.452,80"
$ws.Range("F18").Value = "$ 2.452,80"
$ws.Range("F19").Value = "$ 2.452,80"
$ws.Range("F20").Value = "$ 3.239,58"
$ws.Range("F21").Value = "$ 1.110,65"
$ws.Range("F22").Value = "$ 1.758,58"
$ws.Range("F23").Value = "$ 1.110,65"
$ws.Range("F24").Value = "$ 1.064,37"
$ws.Range("F27").Value = "Sin precio"
$ws.Range("F28").Value = "$ 3.332,13"
$ws.Range("F29").Value = "$ 1.666,02"
$ws.Range("F30").Value = "$ 2.961,88"
$ws.Range("F31").Value = "$ 407,18"
$ws.Range("F32").Value = "$ 490,48"
$ws.Range("F33").Value = "$ 490,48"
$ws.Range("F34").Value = "$ 407,18"
$ws.Range("F35").Value = "$ 832,97"
$ws.Range("F36").Value = "$ 1.851,15"
$ws.Range("F37").Value = "$ 1.851,15"
$ws.Range("F38").Value = "$ 1.601,34"
$ws.Range("F39").Value = "$ 1.388,50"
$ws.Range("F40").Value = "$ 1.489,86"
$ws.Range("F41").Value = "$ 490,48"
$ws.Range("F42").Value = "$ 1.295,77"
$ws.Range("F43").Value = "$ 509,00"
$ws.Range("F44").Value = "$ 1.295,77"
$ws.Range("F45").Value = "$ 1.138,41"
$ws.Range("F46").Value = "$ 490,48"
$ws.Range("F47").Value = "$ 555,28"
$ws.Range("F48").Value = "$ 1.018,08"
$ws.Range("F51").Value = "$ 1.203,22"
$ws.Range("F52").Value = "$ 740,40"
$ws.Range("F53").Value = "$ 2.980,40"
$ws.Range("F54").Value = "$ 1.851,15"
$ws.Range("F55").Value = "$ 740,40"
$ws.Range("F56").Value = "$ 786,69"
$ws.Range("F57").Value = "$ 509,00"
$ws.Range("F58").Value = "$ 944,04"
$ws.Range("F59").Value = "$ 2.221,39"
$ws.Range("F60").Value = "$ 647,84"
$ws.Range("F61").Value = "$ 647,84"
$ws.Range("F62").Value = "$ 1.110,65"
$ws.Range("F65").Value = "$ 481,22"
$ws.Range("F69").Value = "$ 536,77"
$ws.Range("F70").Value = "$ 536,77"
$ws.Range("F72").Value = "$ 1.388,69"
$ws.Range("F73").Value = "$ 1.119,89"
$ws.Range("F74").Value = "$ 1.455,89"
$ws.Range("F75").Value = "$ 1.175,89"
$ws.Range("F76").Value = "$ 1.203,20"
$ws.Range("F77").Value = "$ 786,68"
$ws.Range("F78").Value = "$ 2.036,27"
$ws.Range("F79").Value = "$ 2.545,36"
$ws.Range("F80").Value = "$ 1.480,90"
$ws.Range("F81").Value = "$ 1.156,93"
$ws.Range("F82").Value = "$ 1.156,93"
$ws.Range("F83").Value = "$ 1.156,93"
$ws.Range("F84").Value = "$ 1.156,93"
$ws.Range("F85").Value = "$ 1.156,93"
$ws.Range("F86").Value = "$ 2.776,77"
$ws.Range("F87").Value = "$ 3.100,73"
$ws.Range("F88").Value = "Sin precio"
$ws.Range("F89").Value = "$ 3.100,73"
$ws.Range("F90").Value = "$ 3.100,73"
$ws.Range("F91").Value = "$ 1.110,65"
$ws.Range("F92").Value = "$ 2.776,77"
